$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DSD")
$ws.Activate()

$ws.Range("F4").Value = "CL_COM_GEO_PICT_L123"
$ws.Range("F5").Value = "CL_HH_PRIMARY_ACTIVITY_INDICATORS"

$ws.Columns.Item(6).ColumnWidth = 36.14

$ws.Range("F8").Select()
